$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Latest HO Xliff Generate Date (Overview!G2) and Correspond Handoff Datetime (de-de!H2)
# share the same underlying string value: "2016-08-27 17:01:47" -> "2016-08-27 17:02:43"
$wsOverview.Range("G2").Value = "2016-08-27 17:02:43"
$wsDeDe.Range("H2").Value = "2016-08-27 17:02:43"

# zh-cn!H2 Correspond Handoff Datetime: "2016-08-27 17:01:42" -> "2016-08-27 17:02:39"
$wsZhCn.Range("H2").Value = "2016-08-27 17:02:39"

# zh-cn!K2 Correspond Handback DateTime: "2016-08-27 17:02:11" -> "2016-08-27 17:03:07"
$wsZhCn.Range("K2").Value = "2016-08-27 17:03:07"

# de-de!K2 Correspond Handback DateTime: "2016-08-27 17:02:18" -> "2016-08-27 17:03:14"
$wsDeDe.Range("K2").Value = "2016-08-27 17:03:14"
